$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 361
$ws.Range("I38").Value = 361
$ws.Range("K38").Value = 1083
$ws.Range("M38").Value = -711
$ws.Range("H70").Value = 2994.5
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15540
$ws.Range("H73").Value = 2994.5
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -16872
$ws.Range("H98").Value = 886
$ws.Range("J98").Value = 1003.3333
$ws.Range("L98").Value = 1003.3333
$ws.Range("N98").Value = -3999.3333
$ws.Range("H107").Value = 1023.86365
$ws.Range("I107").Value = 941.8095
$ws.Range("J107").Value = 2747
$ws.Range("K107").Value = 941.8095
$ws.Range("L107").Value = 2747
$ws.Range("M107").Value = 978.1905
$ws.Range("N107").Value = -6587
$ws.Range("H122").Value = 886
$ws.Range("J122").Value = 1003.3333
$ws.Range("L122").Value = 3009.9999
$ws.Range("N122").Value = -7909.9999
$ws.Range("H132").Value = 6134.4736
$ws.Range("I132").Value = 3903.8
$ws.Range("J132").Value = 14499.5
$ws.Range("K132").Value = 11711.4
$ws.Range("L132").Value = 43498.5
$ws.Range("M132").Value = -9181.400000000001
$ws.Range("N132").Value = -48558.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5499
$ws.Range("J45").Value = 5499
$ws.Range("L45").Value = 5499
$ws.Range("N45").Value = -6253
$ws.Range("H110").Value = 1873.3334
$ws.Range("I110").Value = 2133
$ws.Range("K110").Value = 2133
$ws.Range("M110").Value = -88
$ws.Range("H132").Value = 4139.943
$ws.Range("I132").Value = 2134.4138
$ws.Range("J132").Value = 13833.333
$ws.Range("K132").Value = 6403.241399999999
$ws.Range("L132").Value = 41499.999
$ws.Range("M132").Value = -3873.241399999999
$ws.Range("N132").Value = -46559.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 72000
$ws.Range("J35").Value = 72000
$ws.Range("L35").Value = 72000
$ws.Range("N35").Value = -72620
$ws.Range("H86").Value = 21002
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 26503
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 26503
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -28749
$ws.Range("H89").Value = 21002
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 26503
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 132515
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -143747
$ws.Range("H107").Value = 2020.8334
$ws.Range("I107").Value = 1465.1333
$ws.Range("J107").Value = 4799.3335
$ws.Range("K107").Value = 1465.1333
$ws.Range("L107").Value = 4799.3335
$ws.Range("M107").Value = 454.8667
$ws.Range("N107").Value = -8639.333500000001
$ws.Range("H134").Value = 5035.091
$ws.Range("I134").Value = 2840.0833
$ws.Range("J134").Value = 10888.444
$ws.Range("K134").Value = 8520.249899999999
$ws.Range("L134").Value = 32665.332
$ws.Range("M134").Value = -5985.249899999999
$ws.Range("N134").Value = -37735.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 89997.5
$ws.Range("I23").Value = 89997.5
$ws.Range("K23").Value = 89997.5
$ws.Range("M23").Value = -89757.5
$ws.Range("H27").Value = 89997.5
$ws.Range("I27").Value = 89997.5
$ws.Range("K27").Value = 89997.5
$ws.Range("M27").Value = -89805.5
$ws.Range("H86").Value = 5600
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 6200
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 6200
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -8446
$ws.Range("H89").Value = 5600
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 6200
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 31000
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -42232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12224794
$ws.Range("I4").Value = 15715164
$ws.Range("K4").Value = 47145492
$ws.Range("M4").Value = -47145380
$ws.Range("H8").Value = 528.7857
$ws.Range("I8").Value = 528.7857
$ws.Range("K8").Value = 1586.3571
$ws.Range("M8").Value = -1447.3571
$ws.Range("H132").Value = 1946.6364
$ws.Range("I132").Value = 2370.6667
$ws.Range("J132").Value = 1437.8
$ws.Range("K132").Value = 21336.0003
$ws.Range("L132").Value = 12940.2
$ws.Range("M132").Value = -18806.0003
$ws.Range("N132").Value = -18000.2
$ws.Range("H140").Value = 4473.844
$ws.Range("I140").Value = 14734.375
$ws.Range("K140").Value = 44203.125
$ws.Range("M140").Value = -39023.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5581.1665
$ws.Range("I102").Value = 1995
$ws.Range("J102").Value = 7374.25
$ws.Range("K102").Value = 1995
$ws.Range("L102").Value = 7374.25
$ws.Range("M102").Value = -373
$ws.Range("N102").Value = -10618.25
$ws.Range("H126").Value = 4534.615
$ws.Range("I126").Value = 4704.6665
$ws.Range("K126").Value = 14113.9995
$ws.Range("M126").Value = -11643.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2958.818
$ws.Range("J46").Value = 3154.7
$ws.Range("L46").Value = 3154.7
$ws.Range("N46").Value = -3530.7
$ws.Range("H61").Value = 5205.522
$ws.Range("I61").Value = 4184.9375
$ws.Range("J61").Value = 7538.2856
$ws.Range("K61").Value = 4184.9375
$ws.Range("L61").Value = 7538.2856
$ws.Range("M61").Value = -3982.9375
$ws.Range("N61").Value = -7942.2856
$ws.Range("H68").Value = 12133.286
$ws.Range("I68").Value = 9399
$ws.Range("J68").Value = 12589
$ws.Range("K68").Value = 9399
$ws.Range("L68").Value = 12589
$ws.Range("M68").Value = -8650
$ws.Range("N68").Value = -14087
$ws.Range("H71").Value = 12133.286
$ws.Range("I71").Value = 9399
$ws.Range("J71").Value = 12589
$ws.Range("K71").Value = 46995
$ws.Range("L71").Value = 62945
$ws.Range("M71").Value = -43251
$ws.Range("N71").Value = -70433
$ws.Range("H93").Value = 1637.1143
$ws.Range("I93").Value = 1507.762
$ws.Range("K93").Value = 1507.762
$ws.Range("M93").Value = -259.7619999999999
$ws.Range("H113").Value = 5205.522
$ws.Range("I113").Value = 4184.9375
$ws.Range("J113").Value = 7538.2856
$ws.Range("K113").Value = 4184.9375
$ws.Range("L113").Value = 7538.2856
$ws.Range("M113").Value = -2014.9375
$ws.Range("N113").Value = -11878.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 36184.8
$ws.Range("I81").Value = 60696.06
$ws.Range("J81").Value = 4131.615
$ws.Range("K81").Value = 121392.12
$ws.Range("L81").Value = 8263.23
$ws.Range("M81").Value = -120331.12
$ws.Range("N81").Value = -10385.23
$ws.Range("H84").Value = 36184.8
$ws.Range("I84").Value = 60696.06
$ws.Range("J84").Value = 4131.615
$ws.Range("K84").Value = 606960.6
$ws.Range("L84").Value = 41316.14999999999
$ws.Range("M84").Value = -601656.6
$ws.Range("N84").Value = -51924.14999999999
$ws.Range("H107").Value = 2221.5
$ws.Range("I107").Value = 2532
$ws.Range("J107").Value = 1497
$ws.Range("K107").Value = 7596
$ws.Range("L107").Value = 4491
$ws.Range("M107").Value = -5676
$ws.Range("N107").Value = -8331
$ws.Range("H113").Value = 855.41174
$ws.Range("I113").Value = 900.5
$ws.Range("J113").Value = 645
$ws.Range("K113").Value = 2701.5
$ws.Range("L113").Value = 1935
$ws.Range("M113").Value = -531.5
$ws.Range("N113").Value = -6275
$ws.Range("H132").Value = 4734.2
$ws.Range("I132").Value = 3260.2222
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 9780.6666
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -7250.6666
$ws.Range("N132").Value = -59060

Write-Host "Applied Moogle_Profits updates"